$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# HTCS2 (row 26) driver is now completed: TBD -> Done
$ws.Range("B26").Value = "Done"

# New summary block in H8:J10 (Total / Todo / Done + percentage complete)
$ws.Range("H8").Value = "Total"
$ws.Range("I8").Formula = "=COUNTA(B3:B78)"

$ws.Range("H9").Value = "Todo"
$ws.Range("I9").Formula = '=COUNTIF(B3:B78,"TBD")'

$ws.Range("H10").Value = "Done"
$ws.Range("I10").Formula = "=I8-I9"
$ws.Range("J10").Formula = "=I10/I8"

# Style the new labels like the existing "Status" header (bold) and format the numbers
$ws.Range("J10").NumberFormat = "0%"
$ws.Range("I10").NumberFormat = "0"
$ws.Range("H8:H10").Font.Bold = $true

# Extend the Yes/No conditional formatting range used elsewhere to include the new block
$ws.Range("H8:H10").FormatConditions.Add(1, 3, '"Yes"') | Out-Null
$ws.Range("H8:H10").FormatConditions.Add(1, 3, '"No"') | Out-Null

# Column J no longer needs to fit the old "Status" list, shrink it back towards the default width
$ws.Columns.Item(10).ColumnWidth = 8.25

# Restore the view/selection state recorded in the workbook
$ws.Range("L6").Select()
